# Rearranges the per-row data (Fecha, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion, Precio $/Kg, Kg o Unidades) across rows 2-26.
# The mapping below gives, for each destination row, the source row whose
# values should be copied into it (a single-cycle permutation of rows 2..26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 7
    3  = 13
    4  = 24
    5  = 4
    6  = 21
    7  = 20
    8  = 3
    9  = 16
    10 = 2
    11 = 18
    12 = 25
    13 = 5
    14 = 22
    15 = 12
    16 = 11
    17 = 10
    18 = 17
    19 = 8
    20 = 23
    21 = 26
    22 = 6
    23 = 14
    24 = 15
    25 = 9
    26 = 19
}

$cols = @('D', 'I', 'J', 'K', 'L', 'M', 'N', 'P', 'Q')

# Snapshot the original (pre-edit) values for every relevant column/row so
# that writes to one row never clobber data still needed as a source for
# another row later in the loop.
$original = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 26; $r++) {
        $addr = "$col$r"
        $original[$addr] = $ws.Range($addr).Value2
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value2 = $original[$srcAddr]
    }
}
